$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 478
$ws.Range("F8").Value = 358
$ws.Range("F9").Value = 1776
$ws.Range("F11").Value = 1449
$ws.Range("F12").Value = 829
$ws.Range("F13").Value = 357
$ws.Range("F14").Value = 700
$ws.Range("F15").Value = 12968
$ws.Range("F16").Value = 12904
$ws.Range("F17").Value = 967
$ws.Range("F20").Value = 531
$ws.Range("F22").Value = 604
$ws.Range("F23").Value = 2027
$ws.Range("F24").Value = 42
$ws.Range("F26").Value = 13
$ws.Range("F28").Value = 120
$ws.Range("F29").Value = 260
$ws.Range("F30").Value = 702

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 88
$ws.Range("F6").Value = 21

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F6").Value = 478
$ws.Range("F13").Value = 358
$ws.Range("F14").Value = 1776
$ws.Range("F16").Value = 1449
$ws.Range("F17").Value = 829
$ws.Range("F18").Value = 357
$ws.Range("F19").Value = 88
$ws.Range("F20").Value = 700
$ws.Range("F21").Value = 12968
$ws.Range("F22").Value = 12904
$ws.Range("F23").Value = 967
$ws.Range("F26").Value = 531
$ws.Range("F28").Value = 604
$ws.Range("F29").Value = 21
$ws.Range("F31").Value = 2027
$ws.Range("F32").Value = 42
$ws.Range("F34").Value = 13
$ws.Range("F38").Value = 120
$ws.Range("F39").Value = 260
$ws.Range("F40").Value = 702
